$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.724.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.560.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.65%  "

$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.554.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.587"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000277"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.129.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.677.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.64%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.535.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.880"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -13.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "711.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.56%  "

$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.02%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.101"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0478"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.142"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.374.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.319"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0700"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.65%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("E51").Value = "  -0.03%  "
